$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with revised values (in-place edits, no row shift)
$ws.Range("A6").Value = 381
$ws.Range("C6").Value = 11413.92
$ws.Range("D6").Value = 11443.92
$ws.Range("E6").Value = 1.94
$ws.Range("F6").Value = 19.395

$ws.Range("A10").Value = 396
$ws.Range("B10").Value = "8:43 AM"
$ws.Range("C10").Value = 11862.452
$ws.Range("D10").Value = 11892.452
$ws.Range("F10").Value = 11.024

$ws.Range("A11").Value = 399
$ws.Range("B11").Value = "8:45 AM"
$ws.Range("C11").Value = 11968.34
$ws.Range("D11").Value = 11998.34
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 13.115
$ws.Range("G11").Value = 0

$ws.Range("A12").Value = 403
$ws.Range("B12").Value = "8:47 AM"
$ws.Range("C12").Value = 12088.72
$ws.Range("D12").Value = 12118.72
$ws.Range("F12").Value = 32.41

$ws.Range("A13").Value = 407
$ws.Range("C13").Value = 12193.74
$ws.Range("D13").Value = 12223.74
$ws.Range("E13").Value = 1.355
$ws.Range("F13").Value = 58.128333

$ws.Range("A18").Value = 440
$ws.Range("C18").Value = 13177.29
$ws.Range("D18").Value = 13207.29
$ws.Range("F18").Value = 60.2

$ws.Range("A19").Value = 444
$ws.Range("B19").Value = "9:07 AM"
$ws.Range("C19").Value = 13309.11
$ws.Range("D19").Value = 13339.11
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 67.56
$ws.Range("G19").Value = 0

$ws.Range("A20").Value = 448
$ws.Range("B20").Value = "9:09 AM"
$ws.Range("C20").Value = 13418.07
$ws.Range("D20").Value = 13448.07
$ws.Range("F20").Value = 17.015

$ws.Range("A21").Value = 451
$ws.Range("B21").Value = "9:11 AM"
$ws.Range("C21").Value = 13527.81
$ws.Range("D21").Value = 13557.81
$ws.Range("F21").Value = 40.62

$ws.Range("A22").Value = 454
$ws.Range("B22").Value = "9:12 AM"
$ws.Range("C22").Value = 13619.07
$ws.Range("D22").Value = 13649.07
$ws.Range("F22").Value = 49.685

$ws.Range("A24").Value = 464
$ws.Range("C24").Value = 13906.515
$ws.Range("D24").Value = 13936.515
$ws.Range("F24").Value = 11.275

$ws.Range("A25").Value = 467
$ws.Range("B25").Value = "9:18 AM"
$ws.Range("C25").Value = 14001.8
$ws.Range("D25").Value = 14031.8
$ws.Range("F25").Value = 34.39

$ws.Range("A26").Value = 474
$ws.Range("C26").Value = 14204.495
$ws.Range("D26").Value = 14234.495
$ws.Range("F26").Value = 40.96

$ws.Range("A27").Value = 480
$ws.Range("B27").Value = "9:25 AM"
$ws.Range("C27").Value = 14376.725
$ws.Range("D27").Value = 14406.725
$ws.Range("F27").Value = 15.025

$ws.Range("A28").Value = 484
$ws.Range("B28").Value = "9:27 AM"
$ws.Range("C28").Value = 14495.41
$ws.Range("D28").Value = 14525.41
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 65.56999999999999
$ws.Range("G28").Value = 0

$ws.Range("A29").Value = 488
$ws.Range("B29").Value = "9:29 AM"
$ws.Range("C29").Value = 14621.995
$ws.Range("D29").Value = 14651.995
$ws.Range("E29").Value = 1.415
$ws.Range("F29").Value = 113.915

$ws.Range("A30").Value = 497
$ws.Range("C30").Value = 14894.885
$ws.Range("D30").Value = 14924.885
$ws.Range("F30").Value = 30.35

$ws.Range("A31").Value = 501
$ws.Range("B31").Value = "9:35 AM"
$ws.Range("C31").Value = 15000.94
$ws.Range("D31").Value = 15030.94
$ws.Range("F31").Value = 51.39

# Insert a new row before current row 50 (shifting old rows 40-50 down into 41-51),
# then repopulate rows 40-51 with the target values.
$ws.Rows("50").Insert()

$ws.Range("A40").Value = 1035
$ws.Range("B40").Value = "2:02 PM"
$ws.Range("C40").Value = 31029.245
$ws.Range("D40").Value = 31059.245
$ws.Range("E40").Value = 17.475
$ws.Range("F40").Value = 0.545
$ws.Range("G40").Value = 0

$ws.Range("A41").Value = 1043
$ws.Range("B41").Value = "2:06 PM"
$ws.Range("C41").Value = 31284.45
$ws.Range("D41").Value = 31314.45
$ws.Range("E41").Value = 0
$ws.Range("F41").Value = 1.12
$ws.Range("G41").Value = 0

$ws.Range("A42").Value = 1148
$ws.Range("B42").Value = "2:59 PM"
$ws.Range("C42").Value = 34414.095
$ws.Range("D42").Value = 34444.095
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0.8
$ws.Range("G42").Value = 0

$ws.Range("A43").Value = 1383
$ws.Range("B43").Value = "4:56 PM"
$ws.Range("C43").Value = 41472.38
$ws.Range("D43").Value = 41502.38
$ws.Range("E43").Value = 7.295
$ws.Range("F43").Value = 1.325
$ws.Range("G43").Value = 0

$ws.Range("A44").Value = 1389
$ws.Range("B44").Value = "4:59 PM"
$ws.Range("C44").Value = 41640.85
$ws.Range("D44").Value = 41670.85
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 3.15
$ws.Range("G44").Value = 0

$ws.Range("A45").Value = 1397
$ws.Range("B45").Value = "5:04 PM"
$ws.Range("C45").Value = 41905.92
$ws.Range("D45").Value = 41935.92
$ws.Range("E45").Value = 7.78
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = 1

$ws.Range("A46").Value = 1458
$ws.Range("B46").Value = "5:34 PM"
$ws.Range("C46").Value = 43728.35
$ws.Range("D46").Value = 43758.35
$ws.Range("E46").Value = 3.405
$ws.Range("F46").Value = 0.5649999999999999
$ws.Range("G46").Value = 0

$ws.Range("A47").Value = 1547
$ws.Range("B47").Value = "6:19 PM"
$ws.Range("C47").Value = 46409.81
$ws.Range("D47").Value = 46439.81
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 1.36
$ws.Range("G47").Value = 0

$ws.Range("A48").Value = 1551
$ws.Range("B48").Value = "6:20 PM"
$ws.Range("C48").Value = 46502.28
$ws.Range("D48").Value = 46532.28
$ws.Range("E48").Value = 39.27
$ws.Range("F48").Value = 1.01
$ws.Range("G48").Value = 1

$ws.Range("A49").Value = 1558
$ws.Range("B49").Value = "6:24 PM"
$ws.Range("C49").Value = 46717.07
$ws.Range("D49").Value = 46747.07
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0.64
$ws.Range("G49").Value = 0

$ws.Range("A50").Value = 1561
$ws.Range("B50").Value = "6:25 PM"
$ws.Range("C50").Value = 46813.12
$ws.Range("D50").Value = 46843.12
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 1.03
$ws.Range("G50").Value = 0

$ws.Range("A51").Value = 1608
$ws.Range("B51").Value = "6:49 PM"
$ws.Range("C51").Value = 48231.23
$ws.Range("D51").Value = 48261.23
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0.695
$ws.Range("G51").Value = 0
